$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains numeric-looking text (prices with "." as thousands separators).
# Force text format before assigning so values like "1.00" or "0.0000183" are
# preserved exactly instead of being auto-converted to numbers by Excel.
$dCells = @("D2", "D3", "D5", "D6", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D19", "D20", "D21", "D23", "D24", "D26", "D27", "D28", "D29", "D31", "D33", "D34", "D35", "D36", "D38", "D40", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '93.829.34'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').Value = '3.482.68'
$ws.Range('E3').Value = '  +4.90%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '236.03'
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D6').Value = '625.15'
$ws.Range('E6').Value = '  +1.44%  '
$ws.Range('E7').Value = '  +4.73%  '
$ws.Range('E8').Value = '  +3.73%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = '0.995'
$ws.Range('E10').Value = '  +8.43%  '
$ws.Range('D11').Value = '3.481.79'
$ws.Range('E11').Value = '  +4.80%  '
$ws.Range('D12').Value = '43.03'
$ws.Range('E12').Value = '  +3.32%  '
$ws.Range('D13').Value = '0.201'
$ws.Range('E13').Value = '  +5.38%  '
$ws.Range('D14').Value = '6.25'
$ws.Range('E14').Value = '  +5.67%  '
$ws.Range('D15').Value = '4.131.63'
$ws.Range('E15').Value = '  +5.06%  '
$ws.Range('D16').Value = '93.647.31'
$ws.Range('E16').Value = '  +1.38%  '
$ws.Range('E17').Value = '  +4.00%  '
$ws.Range('E18').Value = '  +5.55%  '
$ws.Range('D19').Value = '3.485.01'
$ws.Range('E19').Value = '  +4.98%  '
$ws.Range('D20').Value = '12.59'
$ws.Range('E20').Value = '  +16.01%  '
$ws.Range('D21').Value = '17.85'
$ws.Range('E21').Value = '  +6.17%  '
$ws.Range('E22').Value = '  +11.53%  '
$ws.Range('D23').Value = '519.65'
$ws.Range('E23').Value = '  +7.36%  '
$ws.Range('D24').Value = '3.37'
$ws.Range('E24').Value = '  +4.29%  '
$ws.Range('E25').Value = '  +9.25%  '
$ws.Range('D26').Value = '0.0000183'
$ws.Range('E26').Value = '  +1.42%  '
$ws.Range('D27').Value = '90.72'
$ws.Range('E27').Value = '  +1.64%  '
$ws.Range('D28').Value = '12.22'
$ws.Range('E28').Value = '  +6.57%  '
$ws.Range('D29').Value = '3.669.67'
$ws.Range('E29').Value = '  +4.97%  '
$ws.Range('E30').Value = '  +13.36%  '
$ws.Range('D31').Value = '11.41'
$ws.Range('E31').Value = '  +3.00%  '
$ws.Range('E32').Value = '  +0.28%  '
$ws.Range('D33').Value = '0.138'
$ws.Range('E33').Value = '  +3.47%  '
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('D35').Value = '0.179'
$ws.Range('E35').Value = '  +6.39%  '
$ws.Range('D36').Value = '29.71'
$ws.Range('E36').Value = '  +5.85%  '
$ws.Range('E37').Value = '  +7.47%  '
$ws.Range('D38').Value = '575.09'
$ws.Range('E38').Value = '  +10.98%  '
$ws.Range('E39').Value = '  +7.59%  '
$ws.Range('D40').Value = '7.54'
$ws.Range('E40').Value = '  +4.27%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '0.921'
$ws.Range('E42').Value = '  +5.68%  '
$ws.Range('E43').Value = '  +2.81%  '
$ws.Range('D44').Value = '23.76'
$ws.Range('E44').Value = '  -0.93%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '0.0420'
$ws.Range('E45').Value = '  +6.81%  '
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').Value = '1.69'
$ws.Range('E46').Value = '  +2.41%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '5.49'
$ws.Range('E47').Value = '  +3.46%  '
$ws.Range('B48').Value = 'MantraDAO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D48').Value = '3.55'
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('D49').Value = '2.14'
$ws.Range('E49').Value = '  +3.08%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '8.16'
$ws.Range('E50').Value = '  +3.21%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value = '53.20'
$ws.Range('E51').Value = '  +2.32%  '
